$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.616.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.70%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.491.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.15%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.24%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.596"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.488.64"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.04%  "
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.127"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.440"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.092.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.67%  "
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.26%  "
$ws.Range("E16").Value = "  +4.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.593.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.489.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "392.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.89%  "
$ws.Range("E22").Value = "  +3.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("E25").Value = "  +5.12%  "
$ws.Range("E26").Value = "  +7.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.37%  "
$ws.Range("E28").Value = "  +2.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.94%  "
$ws.Range("E31").Value = "  +7.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.42"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  +10.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.15"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.67%  "
$ws.Range("E38").Value = "  +3.67%  "
$ws.Range("E39").Value = "  +7.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.86%  "
$ws.Range("E41").Value = "  +3.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.51%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.53"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.78%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.79%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.788.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.07%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0314"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.71%  "
$ws.Range("E48").Value = "  +4.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "350.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.85%  "
$ws.Range("E50").Value = "  +6.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +14.92%  "
